# Apply corrected Diebold-Mariano statistics and p-values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -1.138936913881212
$ws.Range("D2").Value = 0.2669798143104525

$ws.Range("C3").Value = -0.5926238491890118
$ws.Range("D3").Value = 0.5594757112472144

$ws.Range("C4").Value = -0.2108856494167026
$ws.Range("D4").Value = 0.8349176317146658

$ws.Range("C5").Value = -0.3594640672952988
$ws.Range("D5").Value = 0.7226753217388806

$ws.Range("C6").Value = 0.6838061076243678
$ws.Range("D6").Value = 0.5012372567111023

$ws.Range("C7").Value = 1.115527375867339
$ws.Range("D7").Value = 0.2766618222240382

$ws.Range("C8").Value = 1.143300595426778
$ws.Range("D8").Value = 0.2652029009680412

$ws.Range("C9").Value = 0.3864741216587823
$ws.Range("D9").Value = 0.7028582073594207

$ws.Range("C10").Value = 0.1734818279474576
$ws.Range("D10").Value = 0.8638590281524323

$ws.Range("C11").Value = -0.1460487573038325
$ws.Range("D11").Value = 0.885212866590209
